$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column D with sequential values 1-15 for rows 1-15 (e.g. inventory id)
for ($i = 1; $i -le 15; $i++) {
    $ws.Cells.Item($i, 4).Value = $i
}

# Update the selected cell to I6 as in the resulting file
$ws.Range("I6").Select()

$wb.Save()
